$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st worksheet) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 15219
$ws1.Range("F3").Value = 19646
$ws1.Range("F5").Value = 191
$ws1.Range("F13").Value = 64
$ws1.Range("F14").Value = 235
$ws1.Range("F17").Value = 1549
$ws1.Range("F20").Value = 128
$ws1.Range("F22").Value = 8332
$ws1.Range("F24").Value = 49
$ws1.Range("F25").Value = 14
$ws1.Range("F26").Value = 74
$ws1.Range("F27").Value = 1289
$ws1.Range("F28").Value = 47
$ws1.Range("F30").Value = 26
$ws1.Range("F31").Value = 6989
$ws1.Range("F34").Value = 196
$ws1.Range("F36").Value = 322
$ws1.Range("F37").Value = 5728
$ws1.Range("F39").Value = 37
$ws1.Range("F40").Value = 33
$ws1.Range("F41").Value = 71

# Sheet "全部类型" (4th worksheet) - column F ("想去人数") updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 15219
$ws4.Range("F3").Value = 19646
$ws4.Range("F5").Value = 191
$ws4.Range("F13").Value = 64
$ws4.Range("F14").Value = 235
$ws4.Range("F16").Value = 0
$ws4.Range("F17").Value = 1549
$ws4.Range("F21").Value = 128
$ws4.Range("F23").Value = 8332
$ws4.Range("F25").Value = 49
$ws4.Range("F26").Value = 14
$ws4.Range("F27").Value = 74
$ws4.Range("F28").Value = 1289
$ws4.Range("F29").Value = 47
$ws4.Range("F31").Value = 26
$ws4.Range("F34").Value = 6989
$ws4.Range("F37").Value = 196
$ws4.Range("F39").Value = 322
$ws4.Range("F40").Value = 5728
$ws4.Range("F42").Value = 37
$ws4.Range("F43").Value = 33
$ws4.Range("F44").Value = 71
